# Flag more RO names
# Inserts 4 new region rows into the detags sheet, shifting subsequent rows down.
# New rows (by their final row number, after all earlier inserts have been applied):
#   row 5  -> abolished
#   row 11 -> Zolochiv
#   row 13 -> Propounded Empathy
#   row 14 -> Crazed Nations CN

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Insert-RegionRow {
    param(
        [int]$Row,
        [string]$Region,
        [string]$Issues,
        [int]$Minor,
        [string]$MinorTimestamp,
        [int]$Major,
        [string]$MajorTimestamp,
        [bool]$NativeEmbassies,
        [string]$Link,
        [string]$Organizations
    )

    # Push the existing row (and everything below it) down by one.
    $ws.Rows.Item($Row).Insert()

    # Match the formatting used by the other "Region" column cells
    # (bold, centered, top-aligned, thin border) by copying it from the
    # row directly above the newly inserted one.
    $ws.Range("A" + ($Row - 1)).Copy()
    $ws.Range("A" + $Row).PasteSpecial(-4122)

    $ws.Range("A" + $Row).Value = $Region
    $ws.Range("B" + $Row).Value = $Issues
    $ws.Range("C" + $Row).Value = $Minor
    $ws.Range("D" + $Row).Value = $MinorTimestamp
    $ws.Range("E" + $Row).Value = $Major
    $ws.Range("F" + $Row).Value = $MajorTimestamp
    $ws.Range("G" + $Row).Value = $NativeEmbassies
    $ws.Range("H" + $Row).Value = $Link
    $ws.Range("I" + $Row).Value = $Organizations
}

Insert-RegionRow 5 "abolished" "RO" 137 "0:02:17" 206 "0:03:26" $false "https://www.nationstates.net/region=abolished" "Unknown"

Insert-RegionRow 11 "Zolochiv" "RO" 421 "0:07:01" 631 "0:10:31" $false "https://www.nationstates.net/region=zolochiv" "Unknown"

Insert-RegionRow 13 "Propounded Empathy" "RO" 733 "0:12:13" 1099 "0:18:19" $false "https://www.nationstates.net/region=propounded_empathy" "Unknown"

Insert-RegionRow 14 "Crazed Nations CN" "RO" 1087 "0:18:07" 1631 "0:27:11" $false "https://www.nationstates.net/region=crazed_nations_cn" "Unknown"
